# Update SWOT Matrix and Zone to Win risk analysis
$d = $word.ActiveDocument

# --- Change 1: Introduction & Methodology intro paragraph -------------------
$d.Content.Find.Execute(
    "Gray Associates Portfolio Analysis evaluates academic programs using a data-driven methodology that plots Market Score (student demand, employment outlook, and competitive positioning) against Program Economics (revenue efficiency and contribution margin). This framework classifies programs into actionable categories" + [char]8212 + "Grow, Sustain, Transform, Evaluate, or Sunset Review" + [char]8212 + "to guide investment and restructuring decisions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Gray Associates Portfolio Analysis evaluates academic programs by plotting Market Score (student demand 40% + employment 40% + competition 20%) against Program Economics (SCH efficiency + cost structure). Programs are classified as Grow, Sustain, Transform, Evaluate, or Sunset Review. Important: FLC does not have a Gray Associates subscription; scores are proxy estimates based on FLC institutional data, not official Gray output.",
    2
)

# --- Change 2: GROW bullet ----------------------------------------------------
$d.Content.Find.Execute(
    "GROW programs (high market + strong economics): Business Admin, Psychology, Engineering, Health Sciences, Computer Info Systems, Exercise Physiology, Accounting.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "GROW programs (high market + strong economics): Business Admin, Psychology, Engineering, Health Sciences, CIS, Exercise Physiology, Accounting show strongest investment case.",
    2
)

# --- Change 3: SUSTAIN bullet -------------------------------------------------
$d.Content.Find.Execute(
    "SUSTAIN programs (solid market, needs efficiency): Environmental programs, Criminology, Biology, Sociology, Teacher Education.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SUSTAIN programs (solid market, needs efficiency): Environmental programs, Criminology, Biology, Sociology, Teacher Education maintain enrollment but need optimization.",
    2
)

# --- Change 4: TRANSFORM bullet -----------------------------------------------
$d.Content.Find.Execute(
    "TRANSFORM programs (weak market, strong economics): English and Mathematics generate revenue but face enrollment pressure - innovate delivery.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "TRANSFORM programs (weak market, strong economics): English and Math generate significant SCH as foundational/service courses " + [char]8212 + " low Market Score reflects major enrollment, not institutional value.",
    2
)

# --- Change 5: EVALUATE/SUNSET bullet -----------------------------------------
$d.Content.Find.Execute(
    "EVALUATE/SUNSET programs (weak market + economics): Political Science, Philosophy, and Art & Design need strategic review for restructuring or phase-out.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "EVALUATE/SUNSET programs: Political Science, Art & Design require strategic review. Note: NAIS is mission-critical and must not be evaluated on enrollment metrics alone.",
    2
)

# --- Change 6: add new disclaimer bullet after the EVALUATE/SUNSET bullet ----
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "EVALUATE/SUNSET programs:*") {
        $target = $para
        break
    }
}
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($i + 1)
$newPara.Range.Text = "Data source disclaimer: FLC does not have a Gray Associates subscription. Scores are proxy estimates based on FLC data, not official Gray Associates output."
